$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1504239.2
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1504239.2
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 4512717.6
$ws.Range("N17").Value = -4513053.6

# Row 69
$ws.Range("H69").Value = 8950
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 8950
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 26850
$ws.Range("N69").Value = -28598

# Row 72
$ws.Range("H72").Value = 8950
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 8950
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 80550
$ws.Range("N72").Value = -89286

# Row 107
$ws.Range("H107").Value = 695.6957
$ws.Range("I107").Value = 772.1579
$ws.Range("J107").Value = 332.5
$ws.Range("K107").Value = 772.1579
$ws.Range("L107").Value = 332.5
$ws.Range("M107").Value = 1147.8421
$ws.Range("N107").Value = -4172.5

# Row 108
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = ""
$ws.Range("N108").Value = 0

# Row 111
$ws.Range("H111").Value = 377.08334
$ws.Range("I111").Value = 292
$ws.Range("J111").Value = 1313
$ws.Range("K111").Value = 876
$ws.Range("L111").Value = 3939
$ws.Range("M111").Value = 2191
$ws.Range("N111").Value = -10073

# Row 112
$ws.Range("H112").Value = 175184.83
$ws.Range("I112").Value = 3999
$ws.Range("J112").Value = 209422
$ws.Range("K112").Value = 11997
$ws.Range("L112").Value = 628266
$ws.Range("M112").Value = -10889
$ws.Range("N112").Value = -630482

# Row 129
$ws.Range("H129").Value = 500007900
$ws.Range("I129").Value = 1000000000
$ws.Range("J129").Value = 15789
$ws.Range("K129").Value = 3000000000
$ws.Range("L129").Value = 47367
$ws.Range("M129").Value = -2999995000
$ws.Range("N129").Value = -57367

# Row 138
$ws.Range("H138").Value = 5686409.5
$ws.Range("I138").Value = 2177.1667
$ws.Range("J138").Value = 6583919.5
$ws.Range("K138").Value = 6531.500100000001
$ws.Range("L138").Value = 19751758.5
$ws.Range("M138").Value = -1391.500100000001
$ws.Range("N138").Value = -19762038.5

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1382.5294
$ws.Range("I2").Value = 730.3
$ws.Range("J2").Value = 2314.2856
$ws.Range("K2").Value = 730.3
$ws.Range("L2").Value = 2314.2856
$ws.Range("M2").Value = -617.3
$ws.Range("N2").Value = -2540.2856

# Row 32
$ws.Range("H32").Value = 13856.072
$ws.Range("I32").Value = 5568.415
$ws.Range("J32").Value = 41308.938
$ws.Range("K32").Value = 5568.415
$ws.Range("L32").Value = 41308.938
$ws.Range("M32").Value = -5281.415
$ws.Range("N32").Value = -41882.938

# Row 61
$ws.Range("H61").Value = 2319.8572
$ws.Range("I61").Value = 1684.742
$ws.Range("J61").Value = 7242
$ws.Range("K61").Value = 1684.742
$ws.Range("L61").Value = 7242
$ws.Range("M61").Value = -1472.742
$ws.Range("N61").Value = -7666

# Row 116
$ws.Range("H116").Value = 1382.5294
$ws.Range("I116").Value = 730.3
$ws.Range("J116").Value = 2314.2856
$ws.Range("K116").Value = 730.3
$ws.Range("L116").Value = 2314.2856
$ws.Range("M116").Value = 1563.7
$ws.Range("N116").Value = -6902.2856

# Row 132
$ws.Range("H132").Value = 2437.1191
$ws.Range("I132").Value = 2193.8206
$ws.Range("J132").Value = 5600
$ws.Range("K132").Value = 6581.4618
$ws.Range("L132").Value = 16800
$ws.Range("M132").Value = -4051.4618
$ws.Range("N132").Value = -21860

# Row 136
$ws.Range("H136").Value = 2319.8572
$ws.Range("I136").Value = 1684.742
$ws.Range("J136").Value = 7242
$ws.Range("K136").Value = 5054.226
$ws.Range("L136").Value = 21726
$ws.Range("M136").Value = -2504.226
$ws.Range("N136").Value = -26826

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1382.5294
$ws.Range("I3").Value = 730.3
$ws.Range("J3").Value = 2314.2856
$ws.Range("K3").Value = 730.3
$ws.Range("L3").Value = 2314.2856
$ws.Range("M3").Value = -616.3
$ws.Range("N3").Value = -2542.2856

# Row 129
$ws.Range("H129").Value = 70000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 70000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 70000
$ws.Range("N129").Value = -80000

# Row 134
$ws.Range("H134").Value = 1271.9
$ws.Range("I134").Value = 1177.8276
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 3533.4828
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -998.4828000000002
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2253.182
$ws.Range("I31").Value = 1832.1765
$ws.Range("J31").Value = 3684.6
$ws.Range("K31").Value = 1832.1765
$ws.Range("L31").Value = 3684.6
$ws.Range("M31").Value = -1537.1765
$ws.Range("N31").Value = -4274.6

# Row 34
$ws.Range("H34").Value = 2253.182
$ws.Range("I34").Value = 1832.1765
$ws.Range("J34").Value = 3684.6
$ws.Range("K34").Value = 1832.1765
$ws.Range("L34").Value = 3684.6
$ws.Range("M34").Value = -1630.1765
$ws.Range("N34").Value = -4088.6

# Row 58
$ws.Range("H58").Value = 2072.5667
$ws.Range("I58").Value = 1608.6957
$ws.Range("J58").Value = 3596.7144
$ws.Range("K58").Value = 1608.6957
$ws.Range("L58").Value = 3596.7144
$ws.Range("M58").Value = -1405.6957
$ws.Range("N58").Value = -4002.7144

# Row 99
$ws.Range("H99").Value = 3126.6667
$ws.Range("I99").Value = 2792
$ws.Range("J99").Value = 3670.5
$ws.Range("K99").Value = 2792
$ws.Range("L99").Value = 3670.5
$ws.Range("M99").Value = -1294
$ws.Range("N99").Value = -6666.5

# Row 126
$ws.Range("H126").Value = 3126.6667
$ws.Range("I126").Value = 2792
$ws.Range("J126").Value = 3670.5
$ws.Range("K126").Value = 8376
$ws.Range("L126").Value = 11011.5
$ws.Range("M126").Value = -5906
$ws.Range("N126").Value = -15951.5

# Row 134
$ws.Range("H134").Value = 9873.784
$ws.Range("I134").Value = 6905.484
$ws.Range("J134").Value = 25210
$ws.Range("K134").Value = 20716.452
$ws.Range("L134").Value = 75630
$ws.Range("M134").Value = -18181.452
$ws.Range("N134").Value = -80700

# Row 136
$ws.Range("H136").Value = 2072.5667
$ws.Range("I136").Value = 1608.6957
$ws.Range("J136").Value = 3596.7144
$ws.Range("K136").Value = 4826.0871
$ws.Range("L136").Value = 10790.1432
$ws.Range("M136").Value = -2276.0871
$ws.Range("N136").Value = -15890.1432

$ws = $wb.Worksheets.Item("CUL")
# Row 60
$ws.Range("H60").Value = 151.6
$ws.Range("I60").Value = 151.6
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 454.8
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -203.8

# Row 107
$ws.Range("H107").Value = 525.5833
$ws.Range("I107").Value = 532.55554
$ws.Range("J107").Value = 504.66666
$ws.Range("K107").Value = 1597.66662
$ws.Range("L107").Value = 1513.99998
$ws.Range("M107").Value = 322.33338
$ws.Range("N107").Value = -5353.999980000001

# Row 112
$ws.Range("H112").Value = 4246.625
$ws.Range("I112").Value = 3992
$ws.Range("J112").Value = 4399.4
$ws.Range("K112").Value = 11976
$ws.Range("L112").Value = 13198.2
$ws.Range("M112").Value = -10868
$ws.Range("N112").Value = -15414.2

# Row 122
$ws.Range("H122").Value = 2324.0417
$ws.Range("I122").Value = 2435.1
$ws.Range("J122").Value = 2244.7144
$ws.Range("K122").Value = 21915.9
$ws.Range("L122").Value = 20202.4296
$ws.Range("M122").Value = -19465.9
$ws.Range("N122").Value = -25102.4296

# Row 138
$ws.Range("H138").Value = 62516500
$ws.Range("I138").Value = 125018000
$ws.Range("J138").Value = 15000
$ws.Range("K138").Value = 375054000
$ws.Range("L138").Value = 45000
$ws.Range("M138").Value = -375048860
$ws.Range("N138").Value = -55280

# Row 140
$ws.Range("H140").Value = 15000
$ws.Range("I140").Value = 15000
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 45000
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -39820

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1624.125
$ws.Range("I122").Value = 1635.1818
$ws.Range("J122").Value = 1599.8
$ws.Range("K122").Value = 4905.5454
$ws.Range("L122").Value = 4799.4
$ws.Range("M122").Value = -2455.5454
$ws.Range("N122").Value = -9699.4

# Row 134
$ws.Range("H134").Value = 29730.4
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 29730.4
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 89191.20000000001
$ws.Range("N134").Value = -94261.20000000001

# Row 136
$ws.Range("H136").Value = 57997.5
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 57997.5
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 173992.5
$ws.Range("N136").Value = -179092.5

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 4920.4287
$ws.Range("I16").Value = 4123.375
$ws.Range("J16").Value = 7471
$ws.Range("K16").Value = 4123.375
$ws.Range("L16").Value = 7471
$ws.Range("M16").Value = -3953.375
$ws.Range("N16").Value = -7811

# Row 111
$ws.Range("H111").Value = 95249.664
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 95249.664
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 95249.664
$ws.Range("N111").Value = -103429.664

# Row 121
$ws.Range("H121").Value = 58899.5
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 58899.5
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 58899.5
$ws.Range("N121").Value = -62393.5

# Row 122
$ws.Range("H122").Value = 6418.1523
$ws.Range("I122").Value = 6115.472
$ws.Range("J122").Value = 7507.8
$ws.Range("K122").Value = 18346.416
$ws.Range("L122").Value = 22523.4
$ws.Range("M122").Value = -15896.416
$ws.Range("N122").Value = -27423.4

# Row 132
$ws.Range("H132").Value = 3387.0571
$ws.Range("I132").Value = 2720.0557
$ws.Range("J132").Value = 4093.2942
$ws.Range("K132").Value = 8160.1671
$ws.Range("L132").Value = 12279.8826
$ws.Range("M132").Value = -5630.1671
$ws.Range("N132").Value = -17339.8826

# Row 136
$ws.Range("H136").Value = 4495.4614
$ws.Range("I136").Value = 3599.15
$ws.Range("J136").Value = 7483.1665
$ws.Range("K136").Value = 10797.45
$ws.Range("L136").Value = 22449.4995
$ws.Range("M136").Value = -8247.450000000001
$ws.Range("N136").Value = -27549.4995

$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 966.6667
$ws.Range("I100").Value = 900
$ws.Range("J100").Value = 1000
$ws.Range("K100").Value = 1800
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1259
$ws.Range("N100").Value = -3082

# Row 121
$ws.Range("H121").Value = 43211.832
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 43211.832
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 43211.832
$ws.Range("N121").Value = -46705.832

# Row 122
$ws.Range("H122").Value = 2176.8918
$ws.Range("I122").Value = 1971.9117
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 5915.7351
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -3465.7351
$ws.Range("N122").Value = -18400

# Row 131
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = ""
$ws.Range("N131").Value = 0

# Row 136
$ws.Range("H136").Value = 1730.7736
$ws.Range("I136").Value = 1531.9131
$ws.Range("J136").Value = 3037.5715
$ws.Range("K136").Value = 4595.7393
$ws.Range("L136").Value = 9112.7145
$ws.Range("M136").Value = -2045.7393
$ws.Range("N136").Value = -14212.7145

# Row 137
$ws.Range("H137").Value = 121249.75
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 121249.75
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 121249.75
$ws.Range("N137").Value = -131449.75
